$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation for 2026/02/10 (火, 17, 201) was recorded, inserted
# right after the existing 2026/02/10 06:00 row (row 803) and before the
# 2026/12/29 block — shifting every following row down by one.
$ws.Rows.Item(804).Insert()

# Leading apostrophe forces the date-look-alike string to stay plain text
# (matches every other date cell in column A, which are inlineStr, not
# real dates).
$ws.Cells.Item(804, 1).Value = "'2026/02/10"
$ws.Cells.Item(804, 2).Value = "火"
$ws.Cells.Item(804, 3).Value = 17
$ws.Cells.Item(804, 4).Value = 201
